$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map (derived from the commit diff).
$updates = [ordered]@{
    'D2' = '27.669.61'
    'E2' = '  +0.12%  '
    'D3' = '1.878.86'
    'E3' = '  +0.94%  '
    'E4' = '  -0.51%  '
    'D5' = '331.56'
    'E5' = '  +2.19%  '
    'D6' = '1.001'
    'E6' = '  -0.65%  '
    'D7' = '0.4731'
    'E7' = '  +5.17%  '
    'D8' = '0.3972'
    'E8' = '  +2.83%  '
    'D9' = '48.22'
    'E9' = '  -2.16%  '
    'D10' = '0.08047'
    'E10' = '  +0.26%  '
    'D11' = '1.025'
    'E11' = '  +0.58%  '
    'D12' = '21.81'
    'E12' = '  +1.48%  '
    'D13' = '1.863.36'
    'E13' = '  -1.53%  '
    'D14' = '5.966'
    'E14' = '  +1.42%  '
    'D15' = '7.177'
    'E15' = '  -0.19%  '
    'D16' = '1.001'
    'E16' = '  -0.81%  '
    'D17' = '87.04'
    'E17' = '  +0.52%  '
    'D18' = '0.00001043'
    'E18' = '  +0.79%  '
    'D19' = '0.06619'
    'E19' = '  +0.87%  '
    'D20' = '17.29'
    'E20' = '  +1.02%  '
    'E21' = '  -0.61%  '
    'D22' = '27.701.11'
    'E22' = '  +0.30%  '
    'D23' = '5.512'
    'E23' = '  +0.00%  '
    'D24' = '11.01'
    'E24' = '  +1.24%  '
    'E25' = '  -0.79%  '
    'D26' = '2.101.61'
    'E26' = '  -0.62%  '
    'D27' = '156.18'
    'E27' = '  +3.04%  '
    'D28' = '20.29'
    'E28' = '  +4.13%  '
    'D29' = '2.094'
    'E29' = '  +2.66%  '
    'D30' = '5.612'
    'E30' = '  +1.26%  '
    'D31' = '122.35'
    'E31' = '  +1.23%  '
    'D32' = '0.9761'
    'E32' = '  +5.17%  '
    'D33' = '0.09574'
    'E33' = '  +1.79%  '
    'D34' = '1.451'
    'E34' = '  -0.67%  '
    'D35' = '3.625'
    'E35' = '  -0.36%  '
    'D36' = '5.317'
    'E36' = '  +0.53%  '
    'D37' = '0.06120'
    'E37' = '  +1.99%  '
    'D38' = '0.02257'
    'E38' = '  +1.07%  '
    'D39' = '1.232'
    'E39' = '  +0.06%  '
    'D40' = '8.141'
    'E40' = '  -3.19%  '
    'D41' = '0.6019'
    'E41' = '  +1.15%  '
    'D42' = '1.001'
    'E42' = '  -0.60%  '
    'E43' = '  +2.57%  '
    'D44' = '10.25'
    'E44' = '  -0.86%  '
    'B45' = 'Decentraland'
    'C45' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'D45' = '0.5721'
    'E45' = '  +1.15%  '
    'B46' = 'WEMIXTOKEN'
    'C46' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'D46' = '1.244'
    'E46' = '  -2.91%  '
    'E47' = '  -2.68%  '
    'D48' = '3.408'
    'E48' = '  -0.29%  '
    'D49' = '1.934'
    'E49' = '  -0.02%  '
    'B50' = 'BabyDogeCoin'
    'C50' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D50' = '0.00000000317'
    'E50' = '  +8.67%  '
    'B51' = 'Cronos'
    'C51' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D51' = '0.06819'
    'E51' = '  -0.72%  '
}

# Matches a plain decimal number (e.g. "331.56", "0.08047") - the shape
# the Price column (D) values take. Anything matching this would be
# silently coerced to a float by Excel (dropping trailing zeros / using
# scientific notation), so those cells are forced to text first.
$numericPattern = "^[0-9]+\.[0-9]+$"

foreach ($cell in $updates.Keys) {
    $value = $updates[$cell]
    $range = $ws.Range($cell)
    if ($cell.StartsWith("D") -and ($value -match $numericPattern)) {
        # Price column: force text so Excel keeps the literal digits/
        # trailing zeros instead of coercing to a float, then restore the
        # default style so no stray number-format style sticks around.
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
